$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-07-2021"
$ws.Range("A76").Style = "Normal"
$ws.Range("B76").Value = -7426
$ws.Range("C76").Value = 344
$ws.Range("D76").Value = -14828
$ws.Range("E76").Value = 263
$ws.Range("F76").Value = -1880
$ws.Range("G76").Value = 8675
